# Update the "想去人数" (want-to-go count) column F values on both the
# "展览" and "全部类型" worksheets, which hold duplicated data tables.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    5  = 163
    7  = 1702
    12 = 1394
    13 = 13
    16 = 52
    20 = 126
    23 = 3262
    24 = 397
    25 = 160
    27 = 12
    29 = 134
    30 = 102
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
